# Auto-generated edit script for cryptos.xlsx price/ranking update
# Applies the coin list re-ranking + updated price/volume figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 '26.279.57'
Set-TextValue 2 5 '  +0.59%  '

# Row 3
Set-TextValue 3 4 '1.663.09'
Set-TextValue 3 5 '  +0.53%  '

# Row 4
Set-TextValue 4 5 '  +0.71%  '

# Row 5
Set-TextValue 5 4 '218.56'
Set-TextValue 5 5 '  +0.13%  '

# Row 6
Set-TextValue 6 4 '0.5316'
Set-TextValue 6 5 '  +0.47%  '

# Row 7
Set-TextValue 7 4 '1.009'

# Row 8
Set-TextValue 8 4 '0.2638'
Set-TextValue 8 5 '  +1.04%  '

# Row 9
Set-TextValue 9 4 '0.06365'
Set-TextValue 9 5 '  +0.48%  '

# Row 10
Set-TextValue 10 4 '20.53'
Set-TextValue 10 5 '  +0.64%  '

# Row 11
Set-TextValue 11 4 '0.07847'
Set-TextValue 11 5 '  +1.16%  '

# Row 12
Set-TextValue 12 2 'WrappedEther'
Set-TextValue 12 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 12 4 '1.702.05'
Set-TextValue 12 5 '  +2.96%  '

# Row 13
Set-TextValue 13 2 'Polkadot'
Set-TextValue 13 3 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 13 4 '4.552'
Set-TextValue 13 5 '  +1.22%  '

# Row 14
Set-TextValue 14 4 '1.892.75'
Set-TextValue 14 5 '  +0.69%  '

# Row 15
Set-TextValue 15 4 '0.5537'
Set-TextValue 15 5 '  +1.41%  '

# Row 16
Set-TextValue 16 4 '0.0₅8195'
Set-TextValue 16 5 '  +0.86%  '

# Row 17
Set-TextValue 17 5 '  +0.63%  '

# Row 18
Set-TextValue 18 2 'Dai'
Set-TextValue 18 3 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 18 4 '1.009'
Set-TextValue 18 5 '  +0.66%  '

# Row 19
Set-TextValue 19 2 'Uniswap'
Set-TextValue 19 3 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 19 4 '4.658'
Set-TextValue 19 5 '  +2.63%  '

# Row 20
Set-TextValue 20 2 'BitcoinCash'
Set-TextValue 20 3 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 20 4 '192.52'
Set-TextValue 20 5 '  -0.46%  '

# Row 21
Set-TextValue 21 2 'Avalanche'
Set-TextValue 21 3 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 21 4 '10.20'
Set-TextValue 21 5 '  +1.58%  '

# Row 22
Set-TextValue 22 2 'Chainlink'
Set-TextValue 22 3 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 22 4 '6.056'
Set-TextValue 22 5 '  +0.96%  '

# Row 23
Set-TextValue 23 2 'BinanceUSD'
Set-TextValue 23 3 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 23 4 '1.011'
Set-TextValue 23 5 '  +0.71%  '

# Row 24
Set-TextValue 24 2 'Monero'
Set-TextValue 24 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 24 4 '145.09'
Set-TextValue 24 5 '  +3.30%  '

# Row 25
Set-TextValue 25 2 'Stellar'
Set-TextValue 25 3 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 25 4 '0.1221'
Set-TextValue 25 5 '  -1.53%  '

# Row 26
Set-TextValue 26 2 'Cosmos'
Set-TextValue 26 3 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 26 4 '7.233'
Set-TextValue 26 5 '  -0.53%  '

# Row 27
Set-TextValue 27 2 'EthereumClassic'
Set-TextValue 27 3 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 27 4 '16.12'
Set-TextValue 27 5 '  -0.21%  '

# Row 28
Set-TextValue 28 2 'Toncoin'
Set-TextValue 28 3 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 28 4 '1.484'
Set-TextValue 28 5 '  +3.56%  '

# Row 29
Set-TextValue 29 2 'Hedera'
Set-TextValue 29 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 29 4 '0.05874'
Set-TextValue 29 5 '  -1.00%  '

# Row 30
Set-TextValue 30 2 'PancakeSwap'
Set-TextValue 30 3 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 30 4 '1.279'
Set-TextValue 30 5 '  +0.32%  '

# Row 31
Set-TextValue 31 2 'InternetComputer(DFINITY)'
Set-TextValue 31 3 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 31 4 '3.589'
Set-TextValue 31 5 '  +2.27%  '

# Row 32
Set-TextValue 32 2 'Filecoin'
Set-TextValue 32 3 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 32 4 '3.304'
Set-TextValue 32 5 '  +2.10%  '

# Row 33
Set-TextValue 33 2 'LidoDAOToken'
Set-TextValue 33 3 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 33 4 '1.616'
Set-TextValue 33 5 '  +4.52%  '

# Row 34
Set-TextValue 34 2 'ARBITRUM'
Set-TextValue 34 3 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 34 4 '0.9590'
Set-TextValue 34 5 '  +1.45%  '

# Row 35
Set-TextValue 35 2 'MXToken'
Set-TextValue 35 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 35 4 '2.814'
Set-TextValue 35 5 '  +1.96%  '

# Row 36
Set-TextValue 36 2 'HuobiToken'
Set-TextValue 36 3 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 36 4 '2.427'
Set-TextValue 36 5 '  +0.58%  '

# Row 37
Set-TextValue 37 2 'ImmutableX'
Set-TextValue 37 3 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 37 4 '0.5811'
Set-TextValue 37 5 '  +3.14%  '

# Row 38
Set-TextValue 38 2 'VeChain'
Set-TextValue 38 3 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 38 4 '0.01612'
Set-TextValue 38 5 '  +0.22%  '

# Row 39
Set-TextValue 39 2 'FraxShare'
Set-TextValue 39 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 39 4 '5.894'
Set-TextValue 39 5 '  +0.88%  '

# Row 40
Set-TextValue 40 2 'TrustWalletToken'
Set-TextValue 40 3 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 40 4 '0.8538'
Set-TextValue 40 5 '  +0.83%  '

# Row 41
Set-TextValue 41 2 'PaxDollar'
Set-TextValue 41 3 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 41 4 '1.009'
Set-TextValue 41 5 '  +0.66%  '

# Row 42
Set-TextValue 42 2 'Maker'
Set-TextValue 42 3 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 42 4 '1.048.11'
Set-TextValue 42 5 '  +3.72%  '

# Row 43
Set-TextValue 43 2 'Quant'
Set-TextValue 43 3 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 43 4 '104.25'
Set-TextValue 43 5 '  +3.12%  '

# Row 44
Set-TextValue 44 2 'RocketPoolETH'
Set-TextValue 44 3 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue 44 4 '1.805.73'
Set-TextValue 44 5 '  +0.54%  '

# Row 45
Set-TextValue 45 2 'Aave'
Set-TextValue 45 3 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 45 4 '57.28'
Set-TextValue 45 5 '  +0.75%  '

# Row 46
Set-TextValue 46 2 'BabyDogeCoin'
Set-TextValue 46 3 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 46 4 '0.0₈108'
Set-TextValue 46 5 '  +2.66%  '

# Row 47
Set-TextValue 47 2 'Frax'
Set-TextValue 47 3 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 47 4 '1.012'
Set-TextValue 47 5 '  +0.91%  '

# Row 48
Set-TextValue 48 2 'Mantle'
Set-TextValue 48 3 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 48 4 '0.4373'
Set-TextValue 48 5 '  +1.97%  '

# Row 49
Set-TextValue 49 2 'EnergySwap'
Set-TextValue 49 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 49 4 '7.949'
Set-TextValue 49 5 '  +2.88%  '

# Row 50
Set-TextValue 50 2 'Cronos'
Set-TextValue 50 3 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 50 4 '0.05163'
Set-TextValue 50 5 '  +0.25%  '

# Row 51
Set-TextValue 51 2 'RenderToken'
Set-TextValue 51 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 51 4 '1.445'
Set-TextValue 51 5 '  -1.58%  '
